# ---------------------------------------------------------------------------
# Adds a "2022-Q1" quarterly fund-holdings sheet to the workbook.
#
# The previous "总计" (summary/roll-up) sheet is renamed to "2022-Q1" and its
# contents are replaced with the new quarter's per-fund breakdown table
# (same column layout as the other quarterly sheets: 基金代码/基金名称/
# 基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
#
# A brand-new "总计" sheet is then appended at the end, holding the rolled-up
# history table (日期/持有数量(只)/持有市值(亿元)) with the new 2022-Q1 row
# added on top of the previously existing rows.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# A sheet that already uses the header/index styling (style index 2) we need
# to replicate onto the new/rewritten sheets - every quarterly sheet uses it.
$styleSource = $wb.Worksheets.Item("2021-Q4")

# --- helper: write one data row of the per-fund breakdown table -----------
# Columns: A index(number) B code(text) C name(text) D scale(text)
#          E position(text) F pct(text) G value(text) H rank(number)
function Set-FundRow($ws, $row, $idx, $code, $name, $scale, $pos, $pct, $val, $rank) {
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = $code
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $scale
    $ws.Cells.Item($row, 5).Value = $pos
    $ws.Cells.Item($row, 6).Value = $pct
    $ws.Cells.Item($row, 7).Value = $val
    $ws.Cells.Item($row, 8).Value = $rank
}

# --- helper: write one data row of the rolled-up history table ------------
# Columns: A index(number) B date(text) C count(number) D value(number)
function Set-TotalRow($ws, $row, $idx, $date, $count, $value) {
    $ws.Cells.Item($row, 1).Value = $idx
    $ws.Cells.Item($row, 2).Value = $date
    $ws.Cells.Item($row, 3).Value = $count
    $ws.Cells.Item($row, 4).Value = $value
}

# ===========================================================================
# Step 1: rename the existing "总计" sheet to "2022-Q1" and repurpose it as
# the new quarter's per-fund holdings sheet.
# ===========================================================================
$ws2022 = $wb.Worksheets.Item("总计")
$ws2022.Name = "2022-Q1"

# Wipe the old roll-up content/formatting entirely before rebuilding.
$ws2022.Cells.Clear()

# Re-create the bordered/bold styling (style index 2 in the original file)
# for the header row and the index column by copying it from a sheet that
# already has it, then overwrite with the real header text + data.
$styleSource.Range("B1:H1").Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)
$styleSource.Range("A2:A14").Copy()
$ws2022.Range("A2:A14").PasteSpecial(-4122)

# Column B (fund code) and columns D/E/F/G (numeric-looking figures) must
# stay text - fund codes can have leading zeros, and the figures need to
# preserve literal formatting like trailing zeros - matching how the other
# quarterly sheets store them. A temporary text NumberFormat forces the
# COM layer to store the literal string instead of auto-converting it to a
# number; ClearFormats afterwards drops the leftover style index again so
# the cells end up with no "s" attribute, same as the target file.
$ws2022.Range("B2:B14").NumberFormat = "@"
$ws2022.Range("D2:G14").NumberFormat = "@"

$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

Set-FundRow $ws2022 2 0 "002408" "中信建投医改灵活配置混合A" "25.68" "94.92" "9.95" "2.5552" 3
Set-FundRow $ws2022 3 1 "001417" "汇添富医疗服务灵活配置混合" "38.45" "77.97" "4.53" "1.7418" 6
Set-FundRow $ws2022 4 2 "015122" "汇添富医疗服务灵活配置混合D" "38.45" "77.97" "4.53" "1.7418" 6
Set-FundRow $ws2022 5 3 "010481" "汇添富高质量成长精选2年持有期混合" "36.73" "72.35" "3.73" "1.3700" 7
Set-FundRow $ws2022 6 4 "007553" "中信建投医改灵活配置混合C" "13.34" "94.92" "9.95" "1.3273" 3
Set-FundRow $ws2022 7 5 "012155" "汇添富成长先锋六个月持有期混合型证券投资基金A" "21.01" "71.54" "4.13" "0.8677" 6
Set-FundRow $ws2022 8 6 "010599" "汇添富高质量成长30一年持有期混合A" "20.34" "74.07" "4.12" "0.8380" 7
Set-FundRow $ws2022 9 7 "010090" "中信建投医药健康混合A" "5.37" "94.87" "9.89" "0.5311" 2
Set-FundRow $ws2022 10 8 "003230" "创金合信医疗保健行业股票A" "8.07" "94.55" "3.95" "0.3188" 10
Set-FundRow $ws2022 11 9 "010091" "中信建投医药健康混合C" "2.25" "94.87" "9.89" "0.2225" 2
Set-FundRow $ws2022 12 10 "003231" "创金合信医疗保健行业股票C" "4.28" "94.55" "3.95" "0.1691" 10
Set-FundRow $ws2022 13 11 "011259" "汇添富高质量成长30一年持有期混合C" "0.90" "74.07" "4.12" "0.0371" 7
Set-FundRow $ws2022 14 12 "012156" "汇添富成长先锋六个月持有期混合型证券投资基金C" "0.57" "71.54" "4.13" "0.0235" 6

# Drop the temporary text-format styling now that the literal strings are
# committed, so these cells carry no "s" attribute (matching the source).
$ws2022.Range("B2:B14").ClearFormats()
$ws2022.Range("D2:G14").ClearFormats()

# ===========================================================================
# Step 2: append a brand-new "总计" sheet (right after "2022-Q1") holding the
# updated roll-up history, with the new 2022-Q1 row on top.
# ===========================================================================
$wsTotal = $wb.Worksheets.Add($null, $ws2022)
$wsTotal.Name = "总计"

$styleSource.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$styleSource.Range("A2:A6").Copy()
$wsTotal.Range("A2:A6").PasteSpecial(-4122)

$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

Set-TotalRow $wsTotal 2 0 "2022-Q1" 13 11.74
Set-TotalRow $wsTotal 3 1 "2021-Q4" 20 11.31
Set-TotalRow $wsTotal 4 2 "2021-Q3" 20 10.95
Set-TotalRow $wsTotal 5 3 "2021-Q2" 14 8.77
Set-TotalRow $wsTotal 6 4 "2020-Q4" 3 0.47
